# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de), the handback step now
# records where the translated content landed: the "Latest Target File"
# (F) and "Latest Handback File" (G) columns get populated (with the same
# kind of hyperlink styling already used for the handoff columns), and the
# "Latest Handback DateTime" (H) column is stamped with the real handback
# time instead of the zero-date placeholder. The Status column also moves
# on from "Ready for handoff" to reflect that handback is complete.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Status text: update every cell that shows it (Overview rolls up both
# language columns, each language sheet shows its own Status column) so
# the shared string is replaced everywhere rather than forked.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- zh-cn sheet: populate Latest Target File (F) / Latest Handback File
# (G) hyperlinks for both data rows, and stamp the handback datetime (H).
$zhcnTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/297cafac3bf6b41d7c14e9f29779c99d9db37dc7/e2e/a.md"
$zhcnHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f43e20546e044337d4e6acfcadf396826efbaf6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcnHandbackFileName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhcnTargetUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnHandbackUrl, "", "", $zhcnHandbackFileName)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhcnTargetUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnHandbackUrl, "", "", $zhcnHandbackFileName)

$zhcn.Range("F2").Style = "HyperLink"
$zhcn.Range("G2").Style = "HyperLink"
$zhcn.Range("F3").Style = "HyperLink"
$zhcn.Range("G3").Style = "HyperLink"

$zhcn.Range("H2").Value = "2016-03-24 12:44:17"
$zhcn.Range("H3").Value = "2016-03-24 12:44:17"

# --- de-de sheet: same shape, its own target/handback URLs and datetime.
$dedeTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/297cafac3bf6b41d7c14e9f29779c99d9db37dc7/e2e/a.md"
$dedeHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac3518ce9413874f12d99f9a04fbefc310ee31cc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dedeHandbackFileName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("F2"), $dedeTargetUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeHandbackUrl, "", "", $dedeHandbackFileName)
$dede.Hyperlinks.Add($dede.Range("F3"), $dedeTargetUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeHandbackUrl, "", "", $dedeHandbackFileName)

$dede.Range("F2").Style = "HyperLink"
$dede.Range("G2").Style = "HyperLink"
$dede.Range("F3").Style = "HyperLink"
$dede.Range("G3").Style = "HyperLink"

$dede.Range("H2").Value = "2016-03-24 12:44:27"
$dede.Range("H3").Value = "2016-03-24 12:44:27"
